# Minor update on sequencing assays
# Applies the logical content changes described by the commit:
#  - dataset_type: "Multiplex Ion Beam Imaging" renamed to "MIBI" and
#    repositioned near the top of the list (alphabetised with the other
#    short acronyms).
#  - preparation_instrument_vendor: added "Leica Biosystems" and "Custom",
#    reordered the existing vendors.
#  - preparation_instrument_model: added "AutoStainer XL" and
#    "Visium CytAssist", reordered the existing models.
#  - .metadata: bumped pav:createdOn timestamp.
#  - main sheet data validation ranges extended to cover the grown lists.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. dataset_type sheet: move "Multiplex Ion Beam Imaging" -> "MIBI"
# ---------------------------------------------------------------------
$dsType = $wb.Worksheets.Item("dataset_type")

# Insert a new row above row 4 (currently "DESI") for the relocated MIBI
# entry, then remove the old entry (now pushed down to row 15).
$dsType.Rows("4").Insert()
$dsType.Range("A4").Value = "MIBI"
$dsType.Range("B4").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172"
$dsType.Rows("15").Delete()

# ---------------------------------------------------------------------
# 2. preparation_instrument_vendor sheet: add Leica Biosystems + Custom,
#    reorder existing vendors.
# ---------------------------------------------------------------------
$vendor = $wb.Worksheets.Item("preparation_instrument_vendor")

$vendor.Range("A1").Value = "Leica Biosystems"
$vendor.Range("B1").Value = "https://identifiers.org/RRID:SCR_023603"
$vendor.Range("A2").Value = "Not applicable"
$vendor.Range("B2").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C48660"
$vendor.Range("A3").Value = "HTX Technologies"
$vendor.Range("B3").Value = "https://identifiers.org/RRID:SCR_023734"
$vendor.Range("A4").Value = "10x Genomics"
$vendor.Range("B4").Value = "https://identifiers.org/RRID:SCR_023672"
$vendor.Range("A5").Value = "Hamamatsu"
$vendor.Range("B5").Value = "https://identifiers.org/RRID:SCR_017105"
$vendor.Range("A6").Value = "Custom"
$vendor.Range("B6").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C126386"
$vendor.Range("A7").Value = "SunChrom"
$vendor.Range("B7").Value = "https://identifiers.org/RRID:SCR_023908"

# ---------------------------------------------------------------------
# 3. preparation_instrument_model sheet: add AutoStainer XL + Visium
#    CytAssist, reorder existing models.
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("preparation_instrument_model")

$model.Range("A1").Value = "NanoZoomer S210"
$model.Range("B1").Value = "https://identifiers.org/RRID:SCR_023760"
$model.Range("A2").Value = "Not applicable"
$model.Range("B2").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C48660"
$model.Range("A3").Value = "Sublimator"
$model.Range("B3").Value = "https://identifiers.org/RRID:SCR_023729"
$model.Range("A4").Value = "Chromium Controller"
$model.Range("B4").Value = "https://identifiers.org/RRID:SCR_019326"
$model.Range("A5").Value = "NanoZoomer S360"
$model.Range("B5").Value = "https://identifiers.org/RRID:SCR_023761"
$model.Range("A6").Value = "NanoZoomer S60"
$model.Range("B6").Value = "https://identifiers.org/RRID:SCR_023762"
$model.Range("A7").Value = "Chromium X"
$model.Range("B7").Value = "https://identifiers.org/RRID:SCR_024537"
$model.Range("A8").Value = "AutoStainer XL"
$model.Range("B8").Value = "https://identifiers.org/RRID:SCR_023957"
$model.Range("A9").Value = "Visium CytAssist"
$model.Range("B9").Value = "https://identifiers.org/RRID:SCR_024570"
$model.Range("A10").Value = "SunCollect Sprayer"
$model.Range("B10").Value = "https://identifiers.org/RRID:SCR_023907"
$model.Range("A11").Value = "M3+ Sprayer"
$model.Range("B11").Value = "https://identifiers.org/RRID:SCR_023731"
$model.Range("A12").Value = "Chromium iX"
$model.Range("B12").Value = "https://identifiers.org/RRID:SCR_024536"
$model.Range("A13").Value = "M5 Sprayer"
$model.Range("B13").Value = "https://identifiers.org/RRID:SCR_023733"
$model.Range("A14").Value = "TM-Sprayer"
$model.Range("B14").Value = "https://identifiers.org/RRID:SCR_023732"

# ---------------------------------------------------------------------
# 4. .metadata sheet: bump pav:createdOn
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-10-20T15:01:07-07:00"

# ---------------------------------------------------------------------
# 5. Main sheet: extend data validation ranges to cover the grown lists.
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("10X Multiome")
$main.Range("H2:H1001").Validation.Delete()
$main.Range("H2:H1001").Validation.Add(3, 1, 1, "=preparation_instrument_vendor!`$A`$1:`$A`$7")
$main.Range("I2:I1001").Validation.Delete()
$main.Range("I2:I1001").Validation.Add(3, 1, 1, "=preparation_instrument_model!`$A`$1:`$A`$14")
